# Atualiza instrução de trabalho
#
# - Unhide the previously-filtered rows and clear the active filter
#   criteria (the "Pendente"/March-2025 filter is removed, the rows it
#   hid become visible again) while keeping the AutoFilter range + sort.
# - Rename "Erick Silva" -> "Erick da Silva" for the rows that referenced
#   him (rows 2-9 in the "ITI" sheet).
# - Move the active selection to D21.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ITI")

# Clear the column filters (Status = "Pendente" and the Data month/year
# group filter) and unhide every row the filter had hidden. This also
# restores sheetPr/autoFilter to "no active filter criteria" while still
# keeping the AutoFilter range and the existing sort state.
$ws.ShowAllData() | Out-Null

# The rows that used to read "Erick Silva" in column B now read
# "Erick da Silva".
for ($r = 2; $r -le 9; $r++) {
    $ws.Cells.Item($r, 2).Value = "Erick da Silva"
}

# Move the selection to D21.
$ws.Activate() | Out-Null
$ws.Range("D21").Select() | Out-Null
